# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 522 on the active sheet,
# shifting the existing rows 522-550 down to 523-551 (dimension grows
# from A1:T550 to A1:T551).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 522.. down by one, copying formatting (incl. the date
# number format on column D) from the row above, just like Excel does
# when a row is inserted through the UI.
$ws.Rows(522).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(522, 1).Value  = 10
$ws.Cells.Item(522, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(522, 3).Value  = "La Araucanía"
$ws.Cells.Item(522, 4).Value  = 44585
$ws.Cells.Item(522, 5).Value  = 9
$ws.Cells.Item(522, 6).Value  = "Fruta"
$ws.Cells.Item(522, 7).Value  = 100104
$ws.Cells.Item(522, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(522, 9).Value  = 100104005
$ws.Cells.Item(522, 10).Value = "Pera"
$ws.Cells.Item(522, 11).Value = "Packham's Triumph"
$ws.Cells.Item(522, 12).Value = "Primera"
$ws.Cells.Item(522, 13).Value = 100
$ws.Cells.Item(522, 14).Value = 15000
$ws.Cells.Item(522, 15).Value = 15000
$ws.Cells.Item(522, 16).Value = 15000
$ws.Cells.Item(522, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(522, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(522, 19).Value = 833
$ws.Cells.Item(522, 20).Value = 18
